# Apply updated dSF (column F) values per repull/recalculation of data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -5
$ws.Range("F6").Value  = -2
$ws.Range("F9").Value  = 3
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = 0
